$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = 20
$ws.Range("B3").Value = 0.5096430919275762

$ws.Range("A4").Value = 30
$ws.Range("B4").Value = 0.5643125960006712

$ws.Range("A5").Value = 40
$ws.Range("B5").Value = 0.5829916928361796
